$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2022-06-24 data update: add new aides counts/amounts (column C = nombre_aides, E = montant_total)
$ws.Range("C13").Value = 187865
$ws.Range("E13").Value = 1168297054

$ws.Range("C48").Value = 150637
$ws.Range("E48").Value = 275745293

$ws.Range("C81").Value = 88359
$ws.Range("E81").Value = 499734008

$ws.Range("C91").Value = 18886
$ws.Range("E91").Value = 75390235

$ws.Range("C121").Value = 1306429
$ws.Range("E121").Value = 2275588183

$ws.Range("C129").Value = 633892
$ws.Range("E129").Value = 3436494215

$ws.Range("C132").Value = 586065
$ws.Range("E132").Value = 3473732208

$ws.Range("C144").Value = 25088
$ws.Range("E144").Value = 92615653

$ws.Range("C186").Value = 236842
$ws.Range("E186").Value = 1190190297

$ws.Range("C207").Value = 154665
$ws.Range("E207").Value = 753774902

$ws.Range("C240").Value = 205941
$ws.Range("E240").Value = 1070112032

$ws.Range("C246").Value = 18838
$ws.Range("E246").Value = 71629583
